$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2363.2727
$ws.Range("I2").Value = 886.5
$ws.Range("J2").Value = 4135.4
$ws.Range("K2").Value = 886.5
$ws.Range("L2").Value = 4135.4
$ws.Range("M2").Value = -773.5
$ws.Range("N2").Value = -4361.4
$ws.Range("H19").Value = 1094.2727
$ws.Range("J19").Value = 1267.25
$ws.Range("L19").Value = 1267.25
$ws.Range("N19").Value = -1617.25
$ws.Range("H97").Value = 2203.8
$ws.Range("J97").Value = 2203.8
$ws.Range("L97").Value = 6611.400000000001
$ws.Range("N97").Value = -7603.400000000001
$ws.Range("H98").Value = 751.53845
$ws.Range("I98").Value = 827.3
$ws.Range("K98").Value = 827.3
$ws.Range("M98").Value = 670.7
$ws.Range("H99").Value = 1137.25
$ws.Range("J99").Value = 1416.3334
$ws.Range("L99").Value = 4249.0002
$ws.Range("N99").Value = -7245.0002
$ws.Range("H106").Value = 17124.273
$ws.Range("I106").Value = 21671.875
$ws.Range("K106").Value = 21671.875
$ws.Range("M106").Value = -21040.875
$ws.Range("H122").Value = 751.53845
$ws.Range("I122").Value = 827.3
$ws.Range("K122").Value = 2481.9
$ws.Range("M122").Value = -31.89999999999964

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2812.4285
$ws.Range("I45").Value = 1942.45
$ws.Range("K45").Value = 1942.45
$ws.Range("M45").Value = -1565.45

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 301.33334
$ws.Range("I5").Value = 301.33334
$ws.Range("K5").Value = 301.33334
$ws.Range("M5").Value = -188.33334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 200
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 200
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -87
$ws.Range("N7").Value = ""
$ws.Range("H86").Value = 4702.8667
$ws.Range("I86").Value = 5020.5
$ws.Range("J86").Value = 3432.3333
$ws.Range("K86").Value = 5020.5
$ws.Range("L86").Value = 3432.3333
$ws.Range("M86").Value = -3897.5
$ws.Range("N86").Value = -5678.3333
$ws.Range("H89").Value = 4702.8667
$ws.Range("I89").Value = 5020.5
$ws.Range("J89").Value = 3432.3333
$ws.Range("K89").Value = 25102.5
$ws.Range("L89").Value = 17161.6665
$ws.Range("M89").Value = -19486.5
$ws.Range("N89").Value = -28393.6665
$ws.Range("H92").Value = 44546
$ws.Range("I92").Value = 44546
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 44546
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -42050
$ws.Range("N92").Value = ""
$ws.Range("H102").Value = 44800
$ws.Range("J102").Value = 44800
$ws.Range("L102").Value = 44800
$ws.Range("N102").Value = -49668
$ws.Range("H107").Value = 240.33333
$ws.Range("I107").Value = 240.33333
$ws.Range("K107").Value = 240.33333
$ws.Range("M107").Value = 1679.66667
$ws.Range("H132").Value = 3131.2703
$ws.Range("I132").Value = 3134.6365
$ws.Range("K132").Value = 9403.9095
$ws.Range("M132").Value = -6873.9095

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 27
$ws.Range("I12").Value = 19.666666
$ws.Range("J12").Value = 30.666666
$ws.Range("K12").Value = 58.999998
$ws.Range("L12").Value = 91.99999800000001
$ws.Range("M12").Value = 114.000002
$ws.Range("N12").Value = -437.999998
$ws.Range("H81").Value = 2342.7144
$ws.Range("I81").Value = 2199.5
$ws.Range("K81").Value = 6598.5
$ws.Range("M81").Value = -5475.5
$ws.Range("H84").Value = 2342.7144
$ws.Range("I84").Value = 2199.5
$ws.Range("K84").Value = 19795.5
$ws.Range("M84").Value = -14179.5
$ws.Range("H92").Value = 1583.3334
$ws.Range("J92").Value = 1625
$ws.Range("L92").Value = 4875
$ws.Range("N92").Value = -7371
$ws.Range("H98").Value = 1581.2
$ws.Range("I98").Value = 2900
$ws.Range("K98").Value = 8700
$ws.Range("M98").Value = -7202
$ws.Range("H122").Value = 1425
$ws.Range("I122").Value = 850
$ws.Range("K122").Value = 7650
$ws.Range("M122").Value = -5200
$ws.Range("H140").Value = 5059.5713
$ws.Range("I140").Value = 1267.4
$ws.Range("K140").Value = 3802.2
$ws.Range("M140").Value = 1377.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 750999.8
$ws.Range("I11").Value = 937500
$ws.Range("J11").Value = 4999
$ws.Range("K11").Value = 937500
$ws.Range("L11").Value = 4999
$ws.Range("M11").Value = -937361
$ws.Range("N11").Value = -5277
$ws.Range("H31").Value = 930.5
$ws.Range("I31").Value = 930.5
$ws.Range("K31").Value = 930.5
$ws.Range("M31").Value = -638.5
$ws.Range("H37").Value = 930.5
$ws.Range("I37").Value = 930.5
$ws.Range("K37").Value = 930.5
$ws.Range("M37").Value = -653.5
$ws.Range("H113").Value = 4599.5
$ws.Range("I113").Value = 2499.3333
$ws.Range("K113").Value = 2499.3333
$ws.Range("M113").Value = -329.3332999999998
$ws.Range("H132").Value = 4228.615
$ws.Range("I132").Value = 2544
$ws.Range("K132").Value = 7632
$ws.Range("M132").Value = -5102
$ws.Range("H136").Value = 32596.334
$ws.Range("J136").Value = 32596.334
$ws.Range("L136").Value = 97789.00199999999
$ws.Range("N136").Value = -102889.002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 894.1429000000001
$ws.Range("I55").Value = 670.0714
$ws.Range("J55").Value = 1342.2858
$ws.Range("K55").Value = 670.0714
$ws.Range("L55").Value = 1342.2858
$ws.Range("M55").Value = -497.0714
$ws.Range("N55").Value = -1688.2858
$ws.Range("H82").Value = 1885.2858
$ws.Range("I82").Value = 1850
$ws.Range("J82").Value = 1932.3334
$ws.Range("K82").Value = 1850
$ws.Range("L82").Value = 1932.3334
$ws.Range("M82").Value = -1489
$ws.Range("N82").Value = -2654.3334
$ws.Range("H85").Value = 1885.2858
$ws.Range("I85").Value = 1850
$ws.Range("J85").Value = 1932.3334
$ws.Range("K85").Value = 1850
$ws.Range("L85").Value = 1932.3334
$ws.Range("M85").Value = -602
$ws.Range("N85").Value = -4428.3334
$ws.Range("H93").Value = 1050
$ws.Range("I93").Value = 2800
$ws.Range("J93").Value = 700
$ws.Range("K93").Value = 2800
$ws.Range("L93").Value = 700
$ws.Range("M93").Value = -1552
$ws.Range("N93").Value = -3196
$ws.Range("H95").Value = 32999.668
$ws.Range("J95").Value = 32999.668
$ws.Range("L95").Value = 32999.668
$ws.Range("N95").Value = -38491.668
$ws.Range("H106").Value = 23000
$ws.Range("J106").Value = 23000
$ws.Range("L106").Value = 23000
$ws.Range("N106").Value = -25524

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2184.4285
$ws.Range("I132").Value = 2184.4285
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6553.2855
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4023.2855
$ws.Range("N132").Value = ""
